# Rename consolidated sheets to their normalized (accent-free, lowercase,
# underscore-separated) names as part of the "dominios de valores" cleanup.

$wb = $excel.ActiveWorkbook

$renames = @{
    1  = "nro_registro_Consolidado"
    2  = "provincia_Consolidado"
    3  = "canton_Consolidado"
    4  = "parroquia_Consolidado"
    5  = "nacionalidad_infractor_Consoli"
    6  = "tipo_transporte_Consolidado"
    7  = "familia_Consolidado"
    8  = "nombre_cientifico_Consolidado"
    9  = "uicn_Consolidado"
    10 = "lista_roja_nacional_Consolidad"
    11 = "cites_Consolidado"
    12 = "sexo_Consolidado"
    13 = "etapa_de_vida_Consolidado"
    14 = "estado_fisico_Consolidado"
    15 = "causal_retencion_Consolidado"
}

foreach ($idx in $renames.Keys) {
    $wb.Worksheets.Item($idx).Name = $renames[$idx]
}
